$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows correspond to sheet rows 2..17 (header row 1 is unchanged).
# Columns A..T: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# then the 16 numeric metric columns (E..T).
$data = @(
    @("FAPs", "Nrp1", "ECs", "FAPs", 1, 0.3333333333333333, 0.2151743333333333, 0.645523, 0.01945888736810517, 0.01945888736810517, 3, 1, 87.94127933333334, 263.823838, 0.4109331243514438, 0.4109331243514437, 18.92270615303045, 170.304355377274, 0.0079963013825783, 0.0079963013825783),
    @("FAPs", "Nrp1", "ECs", "M2", 1, 0.3333333333333333, 0.2151743333333333, 0.645523, 0.01945888736810517, 0.01945888736810517, 3, 1, 52.441971, 157.325913, 0.2450515065683088, 0.2450515065683087, 11.284166148611, 101.557495337499, 0.004768429665697204, 0.004768429665697203),
    @("FAPs", "Nrp1", "ECs", "sCs", 1, 0.3333333333333333, 0.2151743333333333, 0.645523, 0.01945888736810517, 0.01945888736810517, 3, 1, 54.667459, 164.002377, 0.255450795093328, 0.255450795093328, 11.76303404535233, 105.867306408171, 0.004970788249813982, 0.004970788249813981),
    @("FAPs", "Nrp1", "ECs", "Fgf2", 1, 0.3333333333333333, 0.2151743333333333, 0.645523, 0.01945888736810517, 0.01945888736810517, 3, 1, 18.95316166666667, 56.85948500000001, 0.08856457398691947, 0.08856457398691944, 4.07823392618389, 36.704105335655, 0.001723368070015683, 0.001723368070015682),
    @("M2", "Nrp1", "ECs", "FAPs", 3, 1, 8.077181333333334, 24.231544, 0.730444748601188, 0.730444748601188, 3, 1, 87.94127933333334, 263.823838, 0.4109331243514438, 0.4109331243514437, 710.3176598606525, 6392.858938745872, 0.3001639427087911, 0.300163942708791),
    @("M2", "Nrp1", "ECs", "M2", 3, 1, 8.077181333333334, 24.231544, 0.730444748601188, 0.730444748601188, 3, 1, 52.441971, 157.325913, 0.2450515065683088, 0.2450515065683087, 423.583309244408, 3812.249783199672, 0.1789965861096307, 0.1789965861096307),
    @("M2", "Nrp1", "ECs", "sCs", 3, 1, 8.077181333333334, 24.231544, 0.730444748601188, 0.730444748601188, 3, 1, 54.667459, 164.002377, 0.255450795093328, 0.255450795093328, 441.5589793755653, 3974.030814380088, 0.1865926918019196, 0.1865926918019195),
    @("M2", "Nrp1", "ECs", "Fgf2", 3, 1, 8.077181333333334, 24.231544, 0.730444748601188, 0.730444748601188, 3, 1, 18.95316166666667, 56.85948500000001, 0.08856457398691947, 0.08856457398691944, 153.0881236216489, 1377.79311259484, 0.06469152798084671, 0.06469152798084669),
    @("sCs", "Nrp1", "ECs", "FAPs", 3, 1, 0.519176, 1.557528, 0.04695070806875992, 0.04695070806875992, 3, 1, 87.94127933333334, 263.823838, 0.4109331243514438, 0.4109331243514437, 45.65700163916267, 410.9130147524641, 0.01929360115720805, 0.01929360115720805),
    @("sCs", "Nrp1", "ECs", "M2", 3, 1, 0.519176, 1.557528, 0.04695070806875992, 0.04695070806875992, 3, 1, 52.441971, 157.325913, 0.2450515065683088, 0.2450515065683087, 27.226612735896, 245.039514623064, 0.01150534174669847, 0.01150534174669847),
    @("sCs", "Nrp1", "ECs", "sCs", 3, 1, 0.519176, 1.557528, 0.04695070806875992, 0.04695070806875992, 3, 1, 54.667459, 164.002377, 0.255450795093328, 0.255450795093328, 28.382032693784, 255.438294244056, 0.01199359570635945, 0.01199359570635945),
    @("sCs", "Nrp1", "ECs", "Fgf2", 3, 1, 0.519176, 1.557528, 0.04695070806875992, 0.04695070806875992, 3, 1, 18.95316166666667, 56.85948500000001, 0.08856457398691947, 0.08856457398691944, 9.840026661453335, 88.56023995308001, 0.004158169458493945, 0.004158169458493944),
    @("Fgf2", "Nrp1", "ECs", "FAPs", 3, 1, 2.246363333333333, 6.73909, 0.2031456559619469, 0.2031456559619469, 3, 1, 87.94127933333334, 263.823838, 0.4109331243514438, 0.4109331243514437, 197.5480653808245, 1777.93258842742, 0.08347927910286637, 0.08347927910286636),
    @("Fgf2", "Nrp1", "ECs", "M2", 3, 1, 2.246363333333333, 6.73909, 0.2031456559619469, 0.2031456559619469, 3, 1, 52.441971, 157.325913, 0.2450515065683088, 0.2450515065683087, 117.80372078213, 1060.23348703917, 0.04978114904628243, 0.04978114904628243),
    @("Fgf2", "Nrp1", "ECs", "sCs", 3, 1, 2.246363333333333, 6.73909, 0.2031456559619469, 0.2031456559619469, 3, 1, 54.667459, 164.002377, 0.255450795093328, 0.255450795093328, 122.8029754241033, 1105.22677881693, 0.05189371933523502, 0.05189371933523501),
    @("Fgf2", "Nrp1", "ECs", "Fgf2", 3, 1, 2.246363333333333, 6.73909, 0.2031456559619469, 0.2031456559619469, 3, 1, 18.95316166666667, 56.85948500000001, 0.08856457398691947, 0.08856457398691944, 42.5756874187389, 383.18118676865, 0.01799150847756314, 0.01799150847756313)
)

$startRow = 2
for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($startRow + $r, $c + 1).Value = $rowValues[$c]
    }
}
